# Applies the "expansão das análises automáticas" edit:
#  - converts taxa (columns E, F) from fractions to percentages (x100),
#    keeping the existing 0.00% number format so displayed values change accordingly
#  - adds three new columns: apoio_medio (L), contribuicoes (M), media_contribuicoes (N)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the three added columns
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Row 2
$ws.Range("E2").Value = 93.85767790262173
$ws.Range("F2").Value = 61.6121308858739
$ws.Range("L2").Value = 91.42267259885418
$ws.Range("M2").Value = 245680
$ws.Range("N2").Value = 318.2383419689119

# Row 3
$ws.Range("E3").Value = 6.142322097378277
$ws.Range("F3").Value = 70.73170731707317
$ws.Range("L3").Value = 89.66360561945585
$ws.Range("M3").Value = 17873
$ws.Range("N3").Value = 308.1551724137931

# Row 4
$ws.Range("E4").Value = 90.80381471389646
$ws.Range("F4").Value = 93.69842460615155
$ws.Range("L4").Value = 89.12265074751335
$ws.Range("M4").Value = 181999
$ws.Range("N4").Value = 145.7157726180945

# Row 5
$ws.Range("E5").Value = 9.196185286103542
$ws.Range("F5").Value = 99.25925925925925
$ws.Range("L5").Value = 98.94667271041796
$ws.Range("M5").Value = 21647
$ws.Range("N5").Value = 161.544776119403

# Row 6
$ws.Range("E6").Value = 97.80701754385966
$ws.Range("F6").Value = 22.27204783258595
$ws.Range("L6").Value = 19.58294201441566
$ws.Range("M6").Value = 2082
$ws.Range("N6").Value = 13.97315436241611

# Row 7
$ws.Range("E7").Value = 2.192982456140351
$ws.Range("F7").Value = 20
$ws.Range("L7").Value = 19.16882921247149
$ws.Range("M7").Value = 126
$ws.Range("N7").Value = 42
